$wb = $excel.ActiveWorkbook

# summary_statistics
$ws = $wb.Worksheets.Item("summary_statistics")
$ws.Range("B2").Value = 315
$ws.Range("C2").Value = -2.87
$ws.Range("D2").Value = 7.85
$ws.Range("E2").Value = 1.19
$ws.Range("F2").Value = 0.7
$ws.Range("G2").Value = 1.78
$ws.Range("H2").Value = 2.23
$ws.Range("J2").Value = 2.234706616729089
$ws.Range("B3").Value = 315
$ws.Range("G3").Value = 0.32
$ws.Range("B4").Value = 315
$ws.Range("B5").Value = 315
$ws.Range("B6").Value = 315
$ws.Range("E6").Value = 0.12
$ws.Range("G6").Value = 0.33
$ws.Range("B7").Value = 315
$ws.Range("E7").Value = 0.4
$ws.Range("B8").Value = 315
$ws.Range("B9").Value = 315
$ws.Range("E9").Value = 0.32
$ws.Range("G9").Value = 0.47
$ws.Range("B10").Value = 315
$ws.Range("B11").Value = 315
$ws.Range("E11").Value = 0.04
$ws.Range("B12").Value = 315
$ws.Range("B13").Value = 315
$ws.Range("E13").Value = 804.79
$ws.Range("F13").Value = 188.95
$ws.Range("G13").Value = 2319.21
$ws.Range("H13").Value = 502.43
$ws.Range("I13").Value = 26.262
$ws.Range("J13").Value = 528.6904999999999
$ws.Range("B14").Value = 312
$ws.Range("E14").Value = 1.58
$ws.Range("G14").Value = 9.210000000000001
$ws.Range("B15").Value = 312
$ws.Range("E15").Value = 4.08
$ws.Range("F15").Value = 1.15
$ws.Range("G15").Value = 7.79
$ws.Range("B16").Value = 312
$ws.Range("D16").Value = 97.90000000000001
$ws.Range("E16").Value = 13.19
$ws.Range("F16").Value = 3.6
$ws.Range("G16").Value = 21.61
$ws.Range("H16").Value = 14.55
$ws.Range("I16").Value = 0.7000000000000001
$ws.Range("J16").Value = 15.25
$ws.Range("B17").Value = 312
$ws.Range("E17").Value = 24.71
$ws.Range("G17").Value = 25.92
$ws.Range("H17").Value = 27.4
$ws.Range("I17").Value = 6.049999999999999
$ws.Range("J17").Value = 33.45
$ws.Range("B18").Value = 312
$ws.Range("G18").Value = 0.72
$ws.Range("B19").Value = 312
$ws.Range("B20").Value = 312
$ws.Range("E20").Value = 3.88
$ws.Range("G20").Value = 3.1
$ws.Range("B21").Value = 312
$ws.Range("E21").Value = 52.31
$ws.Range("F21").Value = 57.25
$ws.Range("G21").Value = 31.35
$ws.Range("H21").Value = 53.17
$ws.Range("I21").Value = 26.875
$ws.Range("J21").Value = 80.05
$ws.Range("B22").Value = 305
$ws.Range("E22").Value = 50.45
$ws.Range("F22").Value = 49
$ws.Range("G22").Value = 28.26
$ws.Range("H22").Value = 46.5
$ws.Range("I22").Value = 27.5
$ws.Range("J22").Value = 74
$ws.Range("B23").Value = 315
$ws.Range("E23").Value = 51.55
$ws.Range("F23").Value = 53.47
$ws.Range("G23").Value = 27.89
$ws.Range("H23").Value = 48.52
$ws.Range("I23").Value = 26.86
$ws.Range("J23").Value = 75.38

# dichotomous_stats
$ws = $wb.Worksheets.Item("dichotomous_stats")
$ws.Range("B2").Value = 278
$ws.Range("C2").Value = 37
$ws.Range("D2").Value = 0.425
$ws.Range("E2").Value = 1.244
$ws.Range("F2").Value = 0.819
$ws.Range("G2").Value = 1.657
$ws.Range("H2").Value = 0.103
$ws.Range("I2").Value = 53.321
$ws.Range("J2").Value = -0.089
$ws.Range("K2").Value = 0.9399999999999999
$ws.Range("B3").Value = 234
$ws.Range("C3").Value = 81
$ws.Range("D3").Value = 0.18
$ws.Range("E3").Value = 1.24
$ws.Range("F3").Value = 1.06
$ws.Range("G3").Value = 0.8179999999999999
$ws.Range("H3").Value = 0.415
$ws.Range("I3").Value = 150.689
$ws.Range("J3").Value = -0.255
$ws.Range("K3").Value = 0.615
$ws.Range("B4").Value = 232
$ws.Range("C4").Value = 83
$ws.Range("D4").Value = 0.223
$ws.Range("E4").Value = 1.253
$ws.Range("F4").Value = 1.029
$ws.Range("G4").Value = 1.074
$ws.Range("H4").Value = 0.285
$ws.Range("I4").Value = 173.795
$ws.Range("J4").Value = -0.187
$ws.Range("K4").Value = 0.634
$ws.Range("B5").Value = 276
$ws.Range("C5").Value = 39
$ws.Range("D5").Value = 0.458
$ws.Range("E5").Value = 1.251
$ws.Range("F5").Value = 0.793
$ws.Range("G5").Value = 1.415
$ws.Range("H5").Value = 0.163
$ws.Range("I5").Value = 47.55
$ws.Range("J5").Value = -0.193
$ws.Range("K5").Value = 1.108
$ws.Range("B6").Value = 190
$ws.Range("C6").Value = 125
$ws.Range("D6").Value = 0.133
$ws.Range("E6").Value = 1.247
$ws.Range("F6").Value = 1.114
$ws.Range("G6").Value = 0.663
$ws.Range("H6").Value = 0.508
$ws.Range("I6").Value = 286.56
$ws.Range("J6").Value = -0.261
$ws.Range("K6").Value = 0.527
$ws.Range("B7").Value = 225
$ws.Range("C7").Value = 90
$ws.Range("D7").Value = 0.178
$ws.Range("E7").Value = 1.245
$ws.Range("F7").Value = 1.067
$ws.Range("G7").Value = 0.86
$ws.Range("H7").Value = 0.391
$ws.Range("I7").Value = 192.197
$ws.Range("J7").Value = -0.23
$ws.Range("K7").Value = 0.586
$ws.Range("B8").Value = 215
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = 0.198
$ws.Range("E8").Value = 1.257
$ws.Range("F8").Value = 1.059
$ws.Range("G8").Value = 0.973
$ws.Range("H8").Value = 0.331
$ws.Range("I8").Value = 223.945
$ws.Range("J8").Value = -0.203
$ws.Range("K8").Value = 0.598
$ws.Range("B9").Value = 238
$ws.Range("C9").Value = 77
$ws.Range("D9").Value = -0.228
$ws.Range("E9").Value = 1.138
$ws.Range("F9").Value = 1.366
$ws.Range("G9").Value = -1.063
$ws.Range("H9").Value = 0.289
$ws.Range("I9").Value = 150.564
$ws.Range("J9").Value = -0.651
$ws.Range("K9").Value = 0.195
$ws.Range("B10").Value = 301
$ws.Range("C10").Value = 14
$ws.Range("D10").Value = 0.374
$ws.Range("E10").Value = 1.211
$ws.Range("F10").Value = 0.837
$ws.Range("G10").Value = 0.804
$ws.Range("H10").Value = 0.435
$ws.Range("I10").Value = 14.369
$ws.Range("J10").Value = -0.621
$ws.Range("K10").Value = 1.369
$ws.Range("B11").Value = 216
$ws.Range("C11").Value = 99
$ws.Range("D11").Value = 0.295
$ws.Range("E11").Value = 1.287
$ws.Range("F11").Value = 0.992
$ws.Range("G11").Value = 1.486
$ws.Range("H11").Value = 0.139
$ws.Range("I11").Value = 234.422
$ws.Range("J11").Value = -0.096
$ws.Range("K11").Value = 0.6860000000000001

# anovas
$ws = $wb.Worksheets.Item("anovas")
$ws.Range("C2").Value = 118.4011847082842
$ws.Range("D2").Value = 39.46706156942805
$ws.Range("E2").Value = 14.01321851643154
$ws.Range("F2").Value = [double]"1.364254211974835e-08"
$ws.Range("C3").Value = 4.204201248483636
$ws.Range("D3").Value = 1.051050312120909
$ws.Range("E3").Value = 0.3290826750367163
$ws.Range("F3").Value = 0.8583529801424807
$ws.Range("C4").Value = 310.2901090562785
$ws.Range("D4").Value = 16.33105837138308
$ws.Range("E4").Value = 7.043194419655037
$ws.Range("F4").Value = [double]"1.867827595169419e-15"

# continuous_correlations
$ws = $wb.Worksheets.Item("continuous_correlations")
$ws.Range("B2").Value = 0.007
$ws.Range("C2").Value = 0.116
$ws.Range("D2").Value = 0.908
$ws.Range("E2").Value = 313
$ws.Range("F2").Value = -0.104
$ws.Range("G2").Value = 0.117
$ws.Range("B3").Value = -0.047
$ws.Range("C3").Value = -0.832
$ws.Range("D3").Value = 0.406
$ws.Range("E3").Value = 310
$ws.Range("F3").Value = -0.157
$ws.Range("G3").Value = 0.064
$ws.Range("B4").Value = 0.013
$ws.Range("C4").Value = 0.224
$ws.Range("D4").Value = 0.823
$ws.Range("E4").Value = 310
$ws.Range("F4").Value = -0.098
$ws.Range("G4").Value = 0.124
$ws.Range("B5").Value = -0.005
$ws.Range("C5").Value = -0.083
$ws.Range("D5").Value = 0.9340000000000001
$ws.Range("E5").Value = 310
$ws.Range("F5").Value = -0.116
$ws.Range("G5").Value = 0.106
$ws.Range("B6").Value = 0.06900000000000001
$ws.Range("C6").Value = 1.209
$ws.Range("D6").Value = 0.228
$ws.Range("E6").Value = 310
$ws.Range("G6").Value = 0.178
$ws.Range("B7").Value = -0.092
$ws.Range("C7").Value = -1.633
$ws.Range("D7").Value = 0.103
$ws.Range("E7").Value = 310
$ws.Range("F7").Value = -0.201
$ws.Range("G7").Value = 0.019
$ws.Range("B8").Value = -0.09
$ws.Range("C8").Value = -1.588
$ws.Range("D8").Value = 0.113
$ws.Range("E8").Value = 310
$ws.Range("F8").Value = -0.199
$ws.Range("G8").Value = 0.021
$ws.Range("B9").Value = -0.012
$ws.Range("C9").Value = -0.21
$ws.Range("D9").Value = 0.834
$ws.Range("E9").Value = 310
$ws.Range("F9").Value = -0.123
$ws.Range("G9").Value = 0.099
$ws.Range("B10").Value = -0.039
$ws.Range("C10").Value = -0.6889999999999999
$ws.Range("D10").Value = 0.491
$ws.Range("E10").Value = 310
$ws.Range("F10").Value = -0.15
$ws.Range("G10").Value = 0.07199999999999999
$ws.Range("B11").Value = 0.023
$ws.Range("C11").Value = 0.392
$ws.Range("D11").Value = 0.695
$ws.Range("E11").Value = 303
$ws.Range("F11").Value = -0.09
$ws.Range("G11").Value = 0.135
$ws.Range("B12").Value = 0.062
$ws.Range("C12").Value = 1.096
$ws.Range("D12").Value = 0.274
$ws.Range("E12").Value = 313
$ws.Range("F12").Value = -0.049
$ws.Range("G12").Value = 0.171
